$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$dValues = @{
    2  = "71.564.31"
    3  = "3.882.44"
    5  = "605.23"
    6  = "174.97"
    7  = "0.670"
    9  = "0.752"
    10 = "0.178"
    11 = "54.21"
    13 = "11.46"
    14 = "4.504.07"
    15 = "3.885.33"
    16 = "20.98"
    17 = "13.96"
    20 = "71.367.53"
    21 = "440.44"
    22 = "4.78"
    23 = "94.24"
    24 = "3.33"
    25 = "13.91"
    26 = "11.78"
    29 = "10.53"
    30 = "8.79"
    31 = "35.19"
    32 = "13.61"
    34 = "47.95"
    36 = "69.79"
    37 = "632.96"
    38 = "0.437"
    41 = "3.34"
    43 = "3.16"
    46 = "10.25"
    48 = "2.89"
    49 = "2.915.61"
    50 = "0.000279"
    51 = "3.22"
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$row]
}

# --- Column E (Volume 1h) updates ---
$eValues = @{
    2  = "  -1.83%  "
    3  = "  -2.74%  "
    4  = "  +0.02%  "
    5  = "  -2.62%  "
    6  = "  +7.30%  "
    7  = "  -2.47%  "
    8  = "  +0.04%  "
    9  = "  -1.01%  "
    10 = "  +5.44%  "
    11 = "  -0.11%  "
    12 = "  +1.54%  "
    13 = "  +3.34%  "
    14 = "  -2.70%  "
    15 = "  -2.81%  "
    16 = "  +1.36%  "
    17 = "  -1.46%  "
    18 = "  -3.91%  "
    19 = "  -2.19%  "
    20 = "  -1.72%  "
    21 = "  +0.09%  "
    22 = "  -3.45%  "
    23 = "  -2.87%  "
    24 = "  -3.82%  "
    25 = "  -3.76%  "
    26 = "  +3.94%  "
    27 = "  -5.70%  "
    28 = "  -0.02%  "
    29 = "  -0.80%  "
    30 = "  +13.96%  "
    31 = "  -3.30%  "
    32 = "  -2.96%  "
    33 = "  -3.65%  "
    34 = "  -0.49%  "
    35 = "  +11.88%  "
    36 = "  -3.93%  "
    37 = "  -1.11%  "
    38 = "  -0.57%  "
    39 = "  -0.78%  "
    40 = "  +0.17%  "
    41 = "  -1.66%  "
    42 = "  -0.17%  "
    43 = "  +19.12%  "
    46 = "  -3.39%  "
    47 = "  -3.87%  "
    48 = "  -13.95%  "
    49 = "  -0.52%  "
    50 = "  +2.83%  "
    51 = "  -5.78%  "
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}

# --- Rows 44 and 45 swap places (Fetch.AI <-> VeChain) with new values ---
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cellD44 = $ws.Range("D44")
$cellD44.NumberFormat = "@"
$cellD44.Value = "0.0473"
$ws.Range("E44").Value = "  -3.85%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cellD45 = $ws.Range("D45")
$cellD45.NumberFormat = "@"
$cellD45.Value = "2.86"
$ws.Range("E45").Value = "  +7.67%  "
